$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "63.635.63"
$ws.Range("E2").Value = "  -3.21%  "

# Row 3
$ws.Range("D3").Value = "2.607.79"
$ws.Range("E3").Value = "  -2.03%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
Set-TextValue $ws.Range("D5") "571.65"
$ws.Range("E5").Value = "  -4.54%  "

# Row 6
Set-TextValue $ws.Range("D6") "155.22"
$ws.Range("E6").Value = "  -2.87%  "

# Row 7
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.620"
$ws.Range("E8").Value = "  -3.33%  "

# Row 9
$ws.Range("D9").Value = "2.605.89"
$ws.Range("E9").Value = "  -2.02%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.116"
$ws.Range("E10").Value = "  -7.93%  "

# Row 11
Set-TextValue $ws.Range("D11") "5.82"
$ws.Range("E11").Value = "  -0.64%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.379"
$ws.Range("E12").Value = "  -4.84%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.156"
$ws.Range("E13").Value = "  -0.20%  "

# Row 14
Set-TextValue $ws.Range("D14") "27.94"
$ws.Range("E14").Value = "  -4.08%  "

# Row 15
$ws.Range("D15").Value = "3.076.35"
$ws.Range("E15").Value = "  -1.98%  "

# Row 16
Set-TextValue $ws.Range("D16") "0.0000179"
$ws.Range("E16").Value = "  -7.96%  "

# Row 17
$ws.Range("D17").Value = "63.513.18"
$ws.Range("E17").Value = "  -3.24%  "

# Row 18
$ws.Range("D18").Value = "2.617.35"
$ws.Range("E18").Value = "  -2.18%  "

# Row 19
Set-TextValue $ws.Range("D19") "11.93"
$ws.Range("E19").Value = "  -5.02%  "

# Row 20
Set-TextValue $ws.Range("D20") "7.49"
$ws.Range("E20").Value = "  +0.50%  "

# Row 21
Set-TextValue $ws.Range("D21") "4.49"
$ws.Range("E21").Value = "  -6.18%  "

# Row 22
Set-TextValue $ws.Range("D22") "340.36"
$ws.Range("E22").Value = "  -3.93%  "

# Row 23
$ws.Range("E23").Value = "  +0.05%  "

# Row 24
Set-TextValue $ws.Range("D24") "67.23"
$ws.Range("E24").Value = "  -3.76%  "

# Row 25
Set-TextValue $ws.Range("D25") "1.81"
$ws.Range("E25").Value = "  +1.34%  "

# Row 26
Set-TextValue $ws.Range("D26") "0.0000107"
$ws.Range("E26").Value = "  -5.49%  "

# Row 27
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D27") "9.06"
$ws.Range("E27").Value = "  -6.90%  "

# Row 28
$ws.Range("B28").Value = "Bittensor"
$ws.Range("C28").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D28") "578.48"
$ws.Range("E28").Value = "  +2.26%  "

# Row 29
Set-TextValue $ws.Range("D29") "1.55"
$ws.Range("E29").Value = "  -4.21%  "

# Row 30
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D30") "1.00"
$ws.Range("E30").Value = "  -0.02%  "

# Row 31
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D31") "0.160"
$ws.Range("E31").Value = "  -1.97%  "

# Row 32
Set-TextValue $ws.Range("D32") "7.82"
$ws.Range("E32").Value = "  -3.69%  "

# Row 33
Set-TextValue $ws.Range("D33") "2.04"
$ws.Range("E33").Value = "  -4.64%  "

# Row 34
Set-TextValue $ws.Range("D34") "1.72"
$ws.Range("E34").Value = "  -5.86%  "

# Row 35
Set-TextValue $ws.Range("D35") "6.56"
$ws.Range("E35").Value = "  -2.04%  "

# Row 36
Set-TextValue $ws.Range("D36") "5.32"
$ws.Range("E36").Value = "  -3.24%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.400"
$ws.Range("E37").Value = "  -5.28%  "

# Row 38
Set-TextValue $ws.Range("D38") "1.00"
$ws.Range("E38").Value = "  +0.08%  "

# Row 39
Set-TextValue $ws.Range("D39") "19.59"
$ws.Range("E39").Value = "  -4.77%  "

# Row 40
Set-TextValue $ws.Range("D40") "154.10"
$ws.Range("E40").Value = "  -0.03%  "

# Row 41
Set-TextValue $ws.Range("D41") "1.86"
$ws.Range("E41").Value = "  -6.11%  "

# Row 42
$ws.Range("E42").Value = "  -0.04%  "

# Row 43
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D43") "41.53"
$ws.Range("E43").Value = "  -3.10%  "

# Row 44
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D44") "2.48"
$ws.Range("E44").Value = "  -0.69%  "

# Row 45
Set-TextValue $ws.Range("D45") "156.85"
$ws.Range("E45").Value = "  -3.02%  "

# Row 46
Set-TextValue $ws.Range("D46") "23.58"
$ws.Range("E46").Value = "  +0.35%  "

# Row 47
Set-TextValue $ws.Range("D47") "3.85"
$ws.Range("E47").Value = "  -5.79%  "

# Row 48
Set-TextValue $ws.Range("D48") "0.0584"
$ws.Range("E48").Value = "  -5.37%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.629"
$ws.Range("E49").Value = "  -2.54%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.0998"
$ws.Range("E50").Value = "  -1.80%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.0246"
$ws.Range("E51").Value = "  -4.87%  "
